$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to text format so numeric-looking strings (e.g. "1.002", "0.09581")
# are stored verbatim as text instead of being coerced into floating point numbers,
# matching the original inlineStr text cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.209.65'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").Value = '1.871.47'
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '311.60'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("E7").Value = '  -0.90%  '

$ws.Range("E8").Value = '  +1.23%  '

$ws.Range("D9").Value = '0.09581'
$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("E10").Value = '  +3.84%  '

$ws.Range("D11").Value = '40.94'
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").Value = '6.486'
$ws.Range("E12").Value = '  +1.31%  '

$ws.Range("D13").Value = '21.01'
$ws.Range("E13").Value = '  +3.33%  '

$ws.Range("D14").Value = '1.867.43'
$ws.Range("E14").Value = '  +4.35%  '

$ws.Range("D15").Value = '1.002'
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").Value = '7.417'
$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("E18").Value = '  +1.28%  '

$ws.Range("D19").Value = '0.06625'
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").Value = '17.49'
$ws.Range("E20").Value = '  +1.70%  '

$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").Value = '6.161'
$ws.Range("E22").Value = '  +2.06%  '

$ws.Range("D23").Value = '28.271.66'
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("D24").Value = '11.30'
$ws.Range("E24").Value = '  +1.90%  '

$ws.Range("D25").Value = '2.283'
$ws.Range("E25").Value = '  +2.94%  '

$ws.Range("D26").Value = '2.531'
$ws.Range("E26").Value = '  +4.77%  '

$ws.Range("D27").Value = '2.081.69'
$ws.Range("E27").Value = '  +3.68%  '

$ws.Range("D28").Value = '21.20'
$ws.Range("E28").Value = '  +4.19%  '

$ws.Range("D29").Value = '157.65'
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Value = '1.067'
$ws.Range("E31").Value = '  +1.86%  '

$ws.Range("D32").Value = '0.1055'
$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("D33").Value = '5.628'
$ws.Range("E33").Value = '  +0.33%  '

$ws.Range("D34").Value = '3.626'
$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("D35").Value = '0.06751'
$ws.Range("E35").Value = '  -1.62%  '

$ws.Range("D36").Value = '9.554'
$ws.Range("E36").Value = '  +5.86%  '

$ws.Range("D37").Value = '0.02388'
$ws.Range("E37").Value = '  +2.94%  '

$ws.Range("D38").Value = '0.2179'
$ws.Range("E38").Value = '  +0.89%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.6345'
$ws.Range("E39").Value = '  +3.31%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.45'
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("D41").Value = '4.975'
$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("D42").Value = '1.180'
$ws.Range("E42").Value = '  +2.30%  '

$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  +1.75%  '

$ws.Range("D45").Value = '0.6021'
$ws.Range("E45").Value = '  +1.89%  '

$ws.Range("D46").Value = '3.666'
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("D48").Value = '123.82'
$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("D49").Value = '1.987'
$ws.Range("E49").Value = '  +1.89%  '

$ws.Range("D50").Value = '1.195'
$ws.Range("E50").Value = '  +1.46%  '

$ws.Range("D51").Value = '0.06841'
$ws.Range("E51").Value = '  +1.39%  '

# Restore default (Normal) style/format on the price & volume columns so no
# lingering custom number-format style is left attached to these cells.
$ws.Range("D2:E51").Style = "Normal"
